# Updates cryptos list values (price + 1h volume change), and three
# swapped coin-name/link/price/volume rows, per the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.800.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.122.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  -2.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.399"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.117.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.734"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.32%  "

$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.89%  "

$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.505.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.38%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.703.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.107.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.05%  "

$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000205"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "449.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "

$ws.Range("E25").Value = "  -4.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.280.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.137"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.234"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.169"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.165"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.18%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("B38").Value = "MantraDAO"
$ws.Range("C38").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "485.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.439"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.74%  "

$ws.Range("E47").Value = "  +0.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.696"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0338"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.14%  "

$ws.Range("E51").Value = "  -0.24%  "
